$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.395.63'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.844.10'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.29%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9989'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.74%  '

$ws.Range("E6").Value = '  -0.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07533'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.31%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2927'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.48%  '

$ws.Range("E10").Value = '  -0.85%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07715'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.849.36'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.001'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6790'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001040'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.10'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.109.63'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.166'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.427.18'
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '228.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.67%  '

$ws.Range("E21").Value = '  -0.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.455'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.98'
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = '  +0.43%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.374'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.45%  '

$ws.Range("E28").Value = '  -0.45%  '

$ws.Range("E29").Value = '  -0.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.279'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05630'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.100'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.73%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.029'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.840'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.81%  '

$ws.Range("E35").Value = '  +0.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7117'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.71%  '

$ws.Range("E37").Value = '  -0.12%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.245.87'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01808'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.766'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.99%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.316'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9016'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9998'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.88'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.73%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.091'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.52%  '

$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000118'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.3999'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.671'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.892'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.09%  '

$ws.Range("E51").Value = '  -0.40%  '
